$wb = $excel.ActiveWorkbook

# --- Overview sheet: update "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 12:15:06"
$wsOverview.Range("G3").Value = "2016-08-19 12:15:06"

# --- zh-cn sheet: update Priority (E), Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-19 12:14:57"
$wsZhCn.Range("H3").Value = "2016-08-19 12:14:57"
$wsZhCn.Range("K2").Value = "2016-08-19 12:15:29"
$wsZhCn.Range("K3").Value = "2016-08-19 12:15:29"

# --- de-de sheet: update Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 12:15:06"
$wsDeDe.Range("H3").Value = "2016-08-19 12:15:06"
$wsDeDe.Range("K2").Value = "2016-08-19 12:15:36"
$wsDeDe.Range("K3").Value = "2016-08-19 12:15:36"
